# Implement scenario driven by spending inputs
# - Inserts a new leading year column (2014) before the existing 2015 column,
#   shifting the annual spending figures one column to the right.
# - Rescales the annual_constant spending-cap cells (now entered in absolute
#   currency units rather than millions) and gives them scientific formatting.
# - Adds a new "int_perc_treatment_support_relative" program row.
# - Widens column A and moves the active selection to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("spending")

# --- Header row: shift years right by one, insert 2014 at column D ---
$years = @(2014,2015,2016,2017,2018,2019,2020,2021,2022,2023,2024,2025,2026,2027,2028,2029,2030,2031,2032,2033,2034,2035)
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(1, 4 + $i).Value = $years[$i]
}

# --- Row 4 (int_perc_xpert, "specified" spending): shift annual values right
#     by one column, duplicating the old D4 figure into the new D4/E4 pair ---
$row4Values = @(2000000,2000000,1000000,500000,2000000,1000000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000,500000)
for ($i = 0; $i -lt $row4Values.Length; $i++) {
    $ws.Cells.Item(4, 4 + $i).Value = $row4Values[$i]
}

# --- annual_constant rows: bump the spending cap and switch to 0.00E+00 ---
$constantRows = @(2,3,5,6)
foreach ($r in $constantRows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.NumberFormat = "0.00E+00"
    $cell.Value = 10000000000
}

# --- New program row 7: int_perc_treatment_support_relative ---
$ws.Cells.Item(7, 1).Value = "int_perc_treatment_support_relative"
$ws.Cells.Item(7, 2).Value = "annual_constant"
$c7 = $ws.Cells.Item(7, 3)
$c7.NumberFormat = "0.00E+00"
$c7.Value = 10000000000

# --- Column A width ---
$ws.Columns("A").ColumnWidth = 32.833333

# --- Active selection moves to D7 ---
$ws.Range("D7").Select() | Out-Null
